$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$addr, [string]$val) {
    # Force the cell to hold the exact literal text, preventing Excel's
    # COM layer from auto-coercing number-looking strings (e.g. "1.000",
    # "11.30", "0.000006534") into doubles and losing formatting.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '25.708.62'
$ws.Range('E2').Value = '  -3.65%  '
Set-TextValue 'D3' '1.746.07'
$ws.Range('E3').Value = '  -5.68%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '236.74'
$ws.Range('E5').Value = '  -9.92%  '
Set-TextValue 'D6' '1.001'
$ws.Range('E6').Value = '  +0.03%  '
Set-TextValue 'D7' '0.4926'
$ws.Range('E7').Value = '  -8.46%  '
Set-TextValue 'D8' '41.65'
$ws.Range('E8').Value = '  -7.80%  '
Set-TextValue 'D9' '0.2511'
$ws.Range('E9').Value = '  -21.48%  '
Set-TextValue 'D10' '0.06001'
$ws.Range('E10').Value = '  -14.09%  '
Set-TextValue 'D11' '1.744.18'
$ws.Range('E11').Value = '  -5.85%  '
Set-TextValue 'D12' '0.06829'
$ws.Range('E12').Value = '  -12.82%  '
Set-TextValue 'D13' '14.79'
$ws.Range('E13').Value = '  -22.39%  '
Set-TextValue 'D14' '4.455'
$ws.Range('E14').Value = '  -11.90%  '
Set-TextValue 'D15' '76.96'
$ws.Range('E15').Value = '  -14.13%  '
Set-TextValue 'D16' '0.5615'
$ws.Range('E16').Value = '  -27.59%  '
Set-TextValue 'D17' '1.001'
$ws.Range('E17').Value = '  +0.05%  '
Set-TextValue 'D18' '1.001'
$ws.Range('E18').Value = '  +0.05%  '
Set-TextValue 'D19' '25.750.61'
$ws.Range('E19').Value = '  -3.56%  '
Set-TextValue 'D20' '11.30'
$ws.Range('E20').Value = '  -20.33%  '
Set-TextValue 'D21' '0.000006534'
$ws.Range('E21').Value = '  -18.62%  '
Set-TextValue 'D22' '1.965.89'
$ws.Range('E22').Value = '  -5.58%  '
Set-TextValue 'D23' '3.992'
$ws.Range('E23').Value = '  -14.32%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D24' '7.866'
$ws.Range('E24').Value = '  -16.52%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D25' '5.005'
$ws.Range('E25').Value = '  -17.29%  '
Set-TextValue 'D26' '137.28'
$ws.Range('E26').Value = '  -3.83%  '
Set-TextValue 'D27' '1.481'
$ws.Range('E27').Value = '  -12.89%  '
Set-TextValue 'D28' '1.812'
$ws.Range('E28').Value = '  -18.56%  '
$ws.Range('E29').Value = '  -14.46%  '
Set-TextValue 'D30' '101.93'
$ws.Range('E30').Value = '  -8.77%  '
Set-TextValue 'D31' '3.738'
$ws.Range('E31').Value = '  -13.85%  '
Set-TextValue 'D32' '0.07994'
$ws.Range('E32').Value = '  -8.83%  '
Set-TextValue 'D33' '3.370'
$ws.Range('E33').Value = '  -18.35%  '
Set-TextValue 'D34' '0.04394'
$ws.Range('E34').Value = '  -10.00%  '
Set-TextValue 'D35' '1.000'
$ws.Range('E35').Value = '  +0.03%  '
Set-TextValue 'D36' '2.632'
$ws.Range('E36').Value = '  -9.29%  '
Set-TextValue 'D37' '0.9786'
$ws.Range('E37').Value = '  -14.49%  '
Set-TextValue 'D38' '0.6032'
$ws.Range('E38').Value = '  -18.46%  '
Set-TextValue 'D39' '2.670'
$ws.Range('E39').Value = '  -14.22%  '
Set-TextValue 'D40' '1.996'
$ws.Range('E40').Value = '  -15.62%  '
$ws.Range('E41').Value = '  +0.03%  '
Set-TextValue 'D42' '103.05'
$ws.Range('E42').Value = '  -6.02%  '
$ws.Range('E43').Value = '  -14.57%  '
Set-TextValue 'D44' '0.7578'
$ws.Range('E44').Value = '  -16.61%  '
$ws.Range('E45').Value = '  -13.12%  '
$ws.Range('E46').Value = '  -23.47%  '
Set-TextValue 'D47' '0.05246'
$ws.Range('E47').Value = '  -10.16%  '
Set-TextValue 'D48' '0.1065'
$ws.Range('E48').Value = '  -15.20%  '
Set-TextValue 'D49' '30.01'
$ws.Range('E49').Value = '  -14.86%  '
Set-TextValue 'D50' '5.868'
$ws.Range('E50').Value = '  -24.28%  '
Set-TextValue 'D51' '52.26'
$ws.Range('E51').Value = '  -13.72%  '
